$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# Title: "Analysing Changes in Wied Fulija Landfill"
#     -> "Analysing Changes in Maghtab Landfill"
# A single Find/Replace over "Wied Fulija" collapses the old
# spellStart/"Wied"/spellEnd + " " + spellStart/"Fulija"/spellEnd runs
# down to one "Maghtab" run (wrapped in a single proofErr pair), leaving
# the following " Landfill" run untouched.
$d.Content.Find.Execute("Wied Fulija", $false, $false, $false, $false, $false, $true, 1, $false, "Maghtab", 2)

# --- Change 2 -----------------------------------------------------------
# "Study the difference made to the quarry throughout the years: changes
#  in size, layout..." -> same sentence with "quarry" replaced by
# "landfill", split into three runs:
#   "Study the difference made to the " / "landfill" / " throughout the
#   years: changes in size, layout..."
# Locate the word "quarry" (first occurrence in the document, inside that
# sentence). Store the whole-document range in a variable first so
# Find.Execute narrows that same range object down to the hit in place.
$hit = $d.Content
$found = $hit.Find.Execute("quarry", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $qStart = $hit.Start
    $qEnd = $hit.End

    # Replace "quarry" with "landfill" in place (still a single run here).
    $wordRange = $d.Range($qStart, $qEnd)
    $wordRange.Text = "landfill"

    # Force run boundaries around the new word - without changing the
    # resulting formatting - by toggling Bold off then back on over its
    # span. This splits the enclosing run into three runs that all end
    # up sharing identical run properties, matching the target markup.
    $newWordRange = $d.Range($qStart, $qStart + 8)
    $newWordRange.Bold = $false
    $newWordRange.Bold = $true
}
